# Atualização das planilhas de backlog e riscos
# - Reescreve a descrição das ações de mitigação de risco (coluna "Como?")
#   para as linhas "Falta de dedicação dos membros com o projeto" e
#   "Desistência de membros do grupo".
# - Move a seleção ativa para B6 (posição salva pelo autor).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "Realizando uma distribuição das tarefas fazendo com que todos participem"
$ws.Range("F5").Value = "Ensinar todos sobre atividade que todos estão realizando e Transferindo as atividades do desistente para outros membros"

$null = $ws.Range("B6").Select()
